$d = $word.ActiveDocument

$replacements = @(
    @("174×4=696", "297×2=594"),
    @("747×2=1494", "907×5=4535"),
    @("659×4=2636", "401×9=3609"),
    @("693×4=2772", "215×4=860"),
    @("745×7=5215", "130×6=780"),
    @("139×2=278", "906×3=2718"),
    @("951×2=1902", "798×6=4788"),
    @("145×4=580", "316×3=948"),
    @("905×9=8145", "353×6=2118"),
    @("963×4=3852", "980×7=6860"),
    @("600×6=3600", "516×5=2580"),
    @("561×6=3366", "451×4=1804"),
    @("428×2=856", "217×3=651"),
    @("511×4=2044", "559×3=1677"),
    @("329×2=658", "102×9=918"),
    @("813×7=5691", "494×6=2964"),
    @("682×4=2728", "966×2=1932"),
    @("611×2=1222", "103×5=515"),
    @("568×4=2272", "274×7=1918"),
    @("981×5=4905", "921×3=2763"),
    @("171×6=1026", "502×2=1004"),
    @("333×6=1998", "200×4=800"),
    @("743×4=2972", "198×8=1584"),
    @("449×7=3143", "369×4=1476"),
    @("758×6=4548", "318×8=2544")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
